$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 2904
$ws.Range("F3").Value = 21339
$ws.Range("G4").Value = "不可售"
$ws.Range("F5").Value = 3182
$ws.Range("F6").Value = 822
$ws.Range("F8").Value = 531
$ws.Range("F9").Value = 786
$ws.Range("F10").Value = 291
$ws.Range("F13").Value = 130
$ws.Range("F14").Value = 539
$ws.Range("F16").Value = 298
$ws.Range("F18").Value = 437
$ws.Range("F19").Value = 103
$ws.Range("F21").Value = 28
$ws.Range("F22").Value = 49
$ws.Range("F23").Value = 138
$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 125
$ws.Range("F4").Value = 350
$ws.Range("G4").Value = 119
$ws.Range("F13").Value = 162
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 6148
$ws.Range("F5").Value = 1679
$ws.Range("F6").Value = 65
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 6148
$ws.Range("F5").Value = 1679
$ws.Range("F6").Value = 2904
$ws.Range("F7").Value = 21339
$ws.Range("G9").Value = "不可售"
$ws.Range("F10").Value = 125
$ws.Range("F11").Value = 350
$ws.Range("G11").Value = 119
$ws.Range("F12").Value = 3182
$ws.Range("F13").Value = 822
$ws.Range("F15").Value = 65
$ws.Range("F17").Value = 531
$ws.Range("F18").Value = 786
$ws.Range("F19").Value = 291
$ws.Range("F25").Value = 130
$ws.Range("F28").Value = 539
$ws.Range("F32").Value = 298
$ws.Range("F33").Value = 162
$ws.Range("F34").Value = 162
$ws.Range("F36").Value = 437
$ws.Range("F38").Value = 103
$ws.Range("F42").Value = 28
$ws.Range("F43").Value = 49
$ws.Range("F49").Value = 138
